# Add a new "2023" column (T) to the data table, mirroring the style of
# column S for each row, and populate it with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header value for the new year column
$ws.Range("T3").Value = 2023

# Data values for the new year column (rows 4-33)
$values = @(
    52.734251206028382,
    44.646801162600475,
    60.998061560200554,
    41.931627189714625,
    38.177163051511151,
    45.607453560981966,
    50.172884880431361,
    44.112367891063748,
    56.155144351753421,
    37.12775271808399,
    26.579446704517768,
    47.691579663423148,
    49.966474107695483,
    44.339536521432947,
    55.468421253968863,
    45.077411133103766,
    42.341975649266388,
    47.772457765110225,
    46.127136558116561,
    38.861148383596195,
    53.259250196123595,
    77.010722119335071,
    59.189709472566221,
    95.060584781341987,
    61.976853895626128,
    48.390234028455353,
    78.407224173903401,
    39.559787476030614,
    41.059215262778977,
    38.128175110356899
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 4
    $ws.Cells.Item($row, 20).Value = $values[$i]
}

# Copy the style (number format etc.) of column S into the new column T,
# row by row, so the new cells render identically to their neighbours.
for ($row = 3; $row -le 33; $row++) {
    $src = $ws.Cells.Item($row, 19)
    $dst = $ws.Cells.Item($row, 20)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# Reset selection back to the top-left cell so no stale <selection> for T4
# remains in the saved sheetView.
$ws.Range("A1").Select() | Out-Null
